# Apply "Hjemme passive updated meanEMG legmaxROM" edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (meanEMG header row) value updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON): the value that was under D2 moved to C2 with an updated value
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 27.667255443025404

# Row 3 (STR): legmaxROM values in B3 and C3 removed
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Update the selected range shown in the sheet view
$ws.Range("B1:E3").Select()
